$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("exp6(II)")
$ws.Activate()

# --- Fill in newly-measured K=5 column values + extra replicate blocks (rows 27-62) ---
$ws.Range("G27").Value = 1.37638888888888
$ws.Range("G28").Value = 1.39923469387755
$ws.Range("G29").Value = 1.4736111111111101
$ws.Range("G30").Value = 1.4820359281437101
$ws.Range("G31").Value = 1.45757575757575
$ws.Range("G32").Value = 1.51515151515151
$ws.Range("C33").Value = 1.5359477124183001
$ws.Range("D33").Value = 1.10248447204968
$ws.Range("E33").Value = 1
$ws.Range("F33").Value = 1
$ws.Range("G33").Value = 1
$ws.Range("C34").Value = 1.89619883040935
$ws.Range("D34").Value = 1.3859375
$ws.Range("E34").Value = 1.11057692307692
$ws.Range("F34").Value = 1.01687116564417
$ws.Range("G34").Value = 1.0014619883040901
$ws.Range("C35").Value = 2.27678571428571
$ws.Range("D35").Value = 1.5177419354838699
$ws.Range("E35").Value = 1.2439613526569999
$ws.Range("F35").Value = 1.1048951048950999
$ws.Range("G35").Value = 1.03266331658291
$ws.Range("C36").Value = 2.4455882352941098
$ws.Range("D36").Value = 1.76479289940828
$ws.Range("E36").Value = 1.3515625
$ws.Range("F36").Value = 1.1925465838509299
$ws.Range("G36").Value = 1.0796089385474801
$ws.Range("C37").Value = 2.6651234567901199
$ws.Range("D37").Value = 1.8068181818181801
$ws.Range("E37").Value = 1.4644970414201099
$ws.Range("F37").Value = 1.27139037433155
$ws.Range("G37").Value = 1.1451612903225801
$ws.Range("C38").Value = 2.7362804878048701
$ws.Range("D38").Value = 1.9660493827160399
$ws.Range("E38").Value = 1.51344086021505
$ws.Range("F38").Value = 1.3371710526315701
$ws.Range("G38").Value = 1.19938650306748
$ws.Range("C39").Value = 2.8026315789473601
$ws.Range("D39").Value = 2.0454545454545401
$ws.Range("E39").Value = 1.61728395061728
$ws.Range("F39").Value = 1.3801369863013699
$ws.Range("G39").Value = 1.20926966292134
$ws.Range("C40").Value = 2.8159509202453901
$ws.Range("D40").Value = 2.14798850574712
$ws.Range("E40").Value = 1.64786585365853
$ws.Range("F40").Value = 1.4654696132596601
$ws.Range("G40").Value = 1.2671232876712299
$ws.Range("C41").Value = 2.9401197604790399
$ws.Range("D41").Value = 2.2100591715976301
$ws.Range("E41").Value = 1.75552486187845
$ws.Range("F41").Value = 1.4787878787878701
$ws.Range("G41").Value = 1.3139204545454499
$ws.Range("C42").Value = 2.69504310344827
$ws.Range("D42").Value = 2.2361111111111098
$ws.Range("E42").Value = 1.8070652173913
$ws.Range("F42").Value = 1.5536723163841799
$ws.Range("G42").Value = 1.3884180790960401
$ws.Range("C43").Value = 3.1722560975609699
$ws.Range("D43").Value = 2.29458598726114
$ws.Range("E43").Value = 1.8128140703517499
$ws.Range("F43").Value = 1.6150568181818099
$ws.Range("G43").Value = 1.3964497041420101
$ws.Range("C44").Value = 3.2390109890109802
$ws.Range("D44").Value = 2.3294117647058799
$ws.Range("E44").Value = 1.8128415300546401
$ws.Range("F44").Value = 1.6298076923076901
$ws.Range("G44").Value = 1.4388586956521701
$ws.Range("C45").Value = 3.24277456647398
$ws.Range("D45").Value = 2.3818181818181801
$ws.Range("E45").Value = 1.88636363636363
$ws.Range("F45").Value = 1.6640127388535
$ws.Range("G45").Value = 1.4799382716049301
$ws.Range("C46").Value = 3.3963815789473601
$ws.Range("D46").Value = 2.41166666666666
$ws.Range("E46").Value = 1.9294478527607299
$ws.Range("F46").Value = 1.63102409638554
$ws.Range("G46").Value = 1.48479729729729
$ws.Range("C47").Value = 3.32407407407407
$ws.Range("D47").Value = 2.4619205298013198
$ws.Range("E47").Value = 1.94785276073619
$ws.Range("F47").Value = 1.6781767955801099
$ws.Range("G47").Value = 1.51107594936708
$ws.Range("C48").Value = 1.54905063291139
$ws.Range("C49").Value = 1.92901234567901
$ws.Range("C50").Value = 2.2105263157894699
$ws.Range("C51").Value = 2.3734375000000001
$ws.Range("C52").Value = 2.5968208092485501
$ws.Range("C53").Value = 2.6694915254237199
$ws.Range("C54").Value = 2.7704678362573101
$ws.Range("C55").Value = 2.87053571428571
$ws.Range("C56").Value = 2.92837078651685
$ws.Range("C57").Value = 3.0464480874316902
$ws.Range("C58").Value = 3.1967455621301699
$ws.Range("C59").Value = 3.2397959183673399
$ws.Range("C60").Value = 3.2383040935672498
$ws.Range("C61").Value = 3.2515151515151501
$ws.Range("C62").Value = 3.36486486486486

# --- Re-create the chart-backing defined names that Excel regenerates when the
#     box-whisker chart's source data is refreshed/extended (v1.34-.55 duplicate
#     the v1.23-.33 series ranges, written twice as Excel does on each refresh) ---
$n = $wb.Names.Add("_xlchart.v1.34", $ws.Range("B2:B62"))
$n.Visible = $false
$n = $wb.Names.Add("_xlchart.v1.35", $ws.Range("C1"))
$n.Visible = $false
$n = $wb.Names.Add("_xlchart.v1.36", $ws.Range("C2:C62"))
$n.Visible = $false
$n = $wb.Names.Add("_xlchart.v1.37", $ws.Range("D1"))
$n.Visible = $false
$n = $wb.Names.Add("_xlchart.v1.38", $ws.Range("D2:D62"))
$n.Visible = $false
$n = $wb.Names.Add("_xlchart.v1.39", $ws.Range("E1"))
$n.Visible = $false
$n = $wb.Names.Add("_xlchart.v1.40", $ws.Range("E2:E62"))
$n.Visible = $false
$n = $wb.Names.Add("_xlchart.v1.41", $ws.Range("F1"))
$n.Visible = $false
$n = $wb.Names.Add("_xlchart.v1.42", $ws.Range("F2:F62"))
$n.Visible = $false
$n = $wb.Names.Add("_xlchart.v1.43", $ws.Range("G1"))
$n.Visible = $false
$n = $wb.Names.Add("_xlchart.v1.44", $ws.Range("G2:G62"))
$n.Visible = $false
$n = $wb.Names.Add("_xlchart.v1.45", $ws.Range("B2:B62"))
$n.Visible = $false
$n = $wb.Names.Add("_xlchart.v1.46", $ws.Range("C1"))
$n.Visible = $false
$n = $wb.Names.Add("_xlchart.v1.47", $ws.Range("C2:C62"))
$n.Visible = $false
$n = $wb.Names.Add("_xlchart.v1.48", $ws.Range("D1"))
$n.Visible = $false
$n = $wb.Names.Add("_xlchart.v1.49", $ws.Range("D2:D62"))
$n.Visible = $false
$n = $wb.Names.Add("_xlchart.v1.50", $ws.Range("E1"))
$n.Visible = $false
$n = $wb.Names.Add("_xlchart.v1.51", $ws.Range("E2:E62"))
$n.Visible = $false
$n = $wb.Names.Add("_xlchart.v1.52", $ws.Range("F1"))
$n.Visible = $false
$n = $wb.Names.Add("_xlchart.v1.53", $ws.Range("F2:F62"))
$n.Visible = $false
$n = $wb.Names.Add("_xlchart.v1.54", $ws.Range("G1"))
$n.Visible = $false
$n = $wb.Names.Add("_xlchart.v1.55", $ws.Range("G2:G62"))
$n.Visible = $false

# --- Leave the sheet scrolled/selected where the author last left it ---
$ws.Range("J51").Select()
